$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume data (and one new coin insertion: LEO at row 26,
# which shifts the subsequent rows down and drops the final "Aave" row).
$rows = @(
    @{ Row = 2; Coin = "Bitcoin"; Link = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"; Price = "24.442.03"; Volume = "  +10.09%  " },
    @{ Row = 3; Coin = "Ethereum"; Link = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"; Price = "1.678.61"; Volume = "  +5.49%  " },
    @{ Row = 4; Coin = "TetherUSD"; Link = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"; Price = "1.003"; Volume = "  -0.08%  " },
    @{ Row = 5; Coin = "BNB"; Link = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"; Price = "306.00"; Volume = "  +2.59%  " },
    @{ Row = 6; Coin = "USDC"; Link = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"; Price = "0.9971"; Volume = "  +0.65%  " },
    @{ Row = 7; Coin = "XRP"; Link = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"; Price = "0.3683"; Volume = "  +1.53%  " },
    @{ Row = 8; Coin = "Cardano"; Link = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"; Price = "0.3425"; Volume = "  +2.32%  " },
    @{ Row = 9; Coin = "OKB"; Link = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"; Price = "48.05"; Volume = "  +16.43%  " },
    @{ Row = 10; Coin = "Polygon"; Link = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"; Price = "1.161"; Volume = "  +3.84%  " },
    @{ Row = 11; Coin = "Dogecoin"; Link = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"; Price = "0.07216"; Volume = "  +3.86%  " },
    @{ Row = 12; Coin = "BinanceUSD"; Link = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"; Price = "0.9998"; Volume = "  -0.10%  " },
    @{ Row = 13; Coin = "Polkadot"; Link = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"; Price = "6.116"; Volume = "  +5.05%  " },
    @{ Row = 14; Coin = "Solana"; Link = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"; Price = "20.14"; Volume = "  +3.70%  " },
    @{ Row = 15; Coin = "Chainlink"; Link = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; Price = "6.719"; Volume = "  +2.79%  " },
    @{ Row = 16; Coin = "WrappedEther"; Link = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; Price = "1.676.60"; Volume = "  +5.85%  " },
    @{ Row = 17; Coin = "ShibaInu"; Link = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; Price = "0.00001101"; Volume = "  +3.62%  " },
    @{ Row = 18; Coin = "Dai"; Link = "https://coinranking.com/coin/MoTuySvg7+dai-dai"; Price = "0.9970"; Volume = "  +0.60%  " },
    @{ Row = 19; Coin = "TRON"; Link = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; Price = "0.06651"; Volume = "  +1.39%  " },
    @{ Row = 20; Coin = "Litecoin"; Link = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; Price = "80.66"; Volume = "  +6.09%  " },
    @{ Row = 21; Coin = "Avalanche"; Link = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; Price = "16.43"; Volume = "  +3.81%  " },
    @{ Row = 22; Coin = "Uniswap"; Link = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; Price = "6.086"; Volume = "  +2.85%  " },
    @{ Row = 23; Coin = "Cosmos"; Link = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; Price = "12.11"; Volume = "  +4.23%  " },
    @{ Row = 24; Coin = "WrappedBTC"; Link = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; Price = "24.395.61"; Volume = "  +9.94%  " },
    @{ Row = 25; Coin = "Toncoin"; Link = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; Price = "2.419"; Volume = "  +1.77%  " },
    @{ Row = 26; Coin = "LEO"; Link = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"; Price = "3.351"; Volume = "  -3.48%  " },
    @{ Row = 27; Coin = "LidoDAOToken"; Link = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"; Price = "2.652"; Volume = "  +6.06%  " },
    @{ Row = 28; Coin = "Monero"; Link = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; Price = "153.20"; Volume = "  +3.50%  " },
    @{ Row = 29; Coin = "EthereumClassic"; Link = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; Price = "19.41"; Volume = "  +1.34%  " },
    @{ Row = 30; Coin = "WrappedliquidstakedEther2.0"; Link = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; Price = "1.861.97"; Volume = "  +6.14%  " },
    @{ Row = 31; Coin = "BitcoinCash"; Link = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; Price = "127.18"; Volume = "  +4.89%  " },
    @{ Row = 32; Coin = "Filecoin"; Link = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; Price = "6.253"; Volume = "  +6.13%  " },
    @{ Row = 33; Coin = "HuobiToken"; Link = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"; Price = "4.029"; Volume = "  +1.20%  " },
    @{ Row = 34; Coin = "ImmutableX"; Link = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; Price = "0.9751"; Volume = "  +5.99%  " },
    @{ Row = 35; Coin = "Stellar"; Link = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; Price = "0.08430"; Volume = "  +3.64%  " },
    @{ Row = 36; Coin = "WEMIXTOKEN"; Link = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"; Price = "1.692"; Volume = "  +4.83%  " },
    @{ Row = 37; Coin = "Aptos"; Link = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; Price = "12.35"; Volume = "  +5.71%  " },
    @{ Row = 38; Coin = "Hedera"; Link = "https://coinranking.com/coin/jad286TjB+hedera-hbar"; Price = "0.06367"; Volume = "  +6.09%  " },
    @{ Row = 39; Coin = "InternetComputer(DFINITY)"; Link = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; Price = "5.306"; Volume = "  +3.80%  " },
    @{ Row = 40; Coin = "VeChain"; Link = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; Price = "0.02307"; Volume = "  +5.98%  " },
    @{ Row = 41; Coin = "FraxShare"; Link = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; Price = "8.646"; Volume = "  +3.55%  " },
    @{ Row = 42; Coin = "TrustWalletToken"; Link = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"; Price = "1.244"; Volume = "  +0.49%  " },
    @{ Row = 43; Coin = "Algorand"; Link = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"; Price = "0.2088"; Volume = "  +5.32%  " },
    @{ Row = 44; Coin = "TheSandbox"; Link = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"; Price = "0.6089"; Volume = "  +5.34%  " },
    @{ Row = 45; Coin = "Frax"; Link = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"; Price = "0.9971"; Volume = "  +0.50%  " },
    @{ Row = 46; Coin = "PancakeSwap"; Link = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"; Price = "3.760"; Volume = "  -0.15%  " },
    @{ Row = 47; Coin = "EnergySwap"; Link = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; Price = "12.88"; Volume = "  +0.31%  " },
    @{ Row = 48; Coin = "Decentraland"; Link = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"; Price = "0.5872"; Volume = "  +5.71%  " },
    @{ Row = 49; Coin = "Quant"; Link = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"; Price = "125.71"; Volume = "  +0.19%  " },
    @{ Row = 50; Coin = "NEARProtocol"; Link = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; Price = "2.009"; Volume = "  +3.57%  " },
    @{ Row = 51; Coin = "Cronos"; Link = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"; Price = "0.07164"; Volume = "  +6.72%  " }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Range("B$r").Value = $item.Coin
    $ws.Range("C$r").Value = $item.Link
    # Force column D to remain text (matches source data which stores prices as
    # strings such as "24.442.03" or "0.00001101"), otherwise Excel auto-converts
    # the numeric-looking text into a floating point number and loses formatting.
    $ws.Range("D$r").NumberFormat = "@"
    $ws.Range("D$r").Value = $item.Price
    $ws.Range("D$r").Style = "Normal"
    $ws.Range("E$r").Value = $item.Volume
}
